# Generate Report for Handback
# Adds a new handback record for file "bda2774e-0609-4935-b45a-8dadf004bda4.md"
# as row 4 on the Overview, zh-cn and de-de sheets, mirroring the structure that
# was used for the previously added "63d8a547-a403-4136-a001-3dc178b0b903.md" entry.

$wb = $excel.ActiveWorkbook

# Color used by the workbook's custom "HyperLink" cell style (RGB FF6495ED),
# expressed as the BGR integer that the Font.Color property expects.
$hyperlinkColor = 15570276
$dateFormat = "yyyy-mm-dd HH:mm:ss"

$guid = "bda2774e-0609-4935-b45a-8dadf004bda4"
$srcFile = "$guid.md"
$pathAndName = "e2e\$guid.md"
$zhXlf = "$guid.969a135b89d14a6b93bf20c5a47f46a1f90dd106.zh-cn.xlf"
$deXlf = "$guid.969a135b89d14a6b93bf20c5a47f46a1f90dd106.de-de.xlf"

$zhHandoffDate = "2016-09-07 04:58:01"
$zhHandbackDate = "2016-09-07 04:58:33"
$deHandoffDate = "2016-09-07 04:58:14"
$deHandbackDate = "2016-09-07 04:58:41"
# Date shown on the Overview sheet for this row (latest HO xliff generate date)
$overviewDate = "2016-09-07 04:58:14"

$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3) -> add row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3:G3").Copy($wsOverview.Range("A4:G4"))

$wsOverview.Range("A4").Value = $srcFile
$wsOverview.Range("B4").Value = $pathAndName

$hlOverview = $wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$srcFile", "", "", $pathAndName)
$wsOverview.Range("B4").Font.Underline = $true
$wsOverview.Range("B4").Font.Color = $hyperlinkColor

$wsOverview.Range("G4").Value = $overviewDate
$wsOverview.Range("G4").NumberFormat = $dateFormat

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1) -> add row 4
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3:P3").Copy($wsZhCn.Range("A4:P4"))

$wsZhCn.Range("A4").Value = $srcFile
$hlZhA = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0000000000000000000000000000000000000000/e2e/$srcFile", "", "", $srcFile)
$wsZhCn.Range("A4").Font.Underline = $true
$wsZhCn.Range("A4").Font.Color = $hyperlinkColor

$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").Value = $zhHandoffDate
$wsZhCn.Range("H4").NumberFormat = $dateFormat

$wsZhCn.Range("I4").Value = $srcFile
$hlZhI = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$srcFile", "", "", $srcFile)
$wsZhCn.Range("I4").Font.Underline = $true
$wsZhCn.Range("I4").Font.Color = $hyperlinkColor

$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").Value = $zhHandbackDate
$wsZhCn.Range("K4").NumberFormat = $dateFormat

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2) -> add row 4
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3:P3").Copy($wsDeDe.Range("A4:P4"))

$wsDeDe.Range("A4").Value = $srcFile
$hlDeA = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0000000000000000000000000000000000000000/e2e/$srcFile", "", "", $srcFile)
$wsDeDe.Range("A4").Font.Underline = $true
$wsDeDe.Range("A4").Font.Color = $hyperlinkColor

$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").Value = $deHandoffDate
$wsDeDe.Range("H4").NumberFormat = $dateFormat

$wsDeDe.Range("I4").Value = $srcFile
$hlDeI = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/$srcFile", "", "", $srcFile)
$wsDeDe.Range("I4").Font.Underline = $true
$wsDeDe.Range("I4").Font.Color = $hyperlinkColor

$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").Value = $deHandbackDate
$wsDeDe.Range("K4").NumberFormat = $dateFormat

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Host "Added handback row for $srcFile to Overview, zh-cn and de-de sheets."
